# Apply updated crypto price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.278.65"
$ws.Range("E2").Value = "  +5.68%  "
$ws.Range("D3").Value = "3.004.06"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.60"
$ws.Range("E5").Value = "  +2.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.94"
$ws.Range("E6").Value = "  +12.86%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("D9").Value = "2.999.78"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.51"
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("E11").Value = "  +3.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +5.28%  "
$ws.Range("E13").Value = "  +5.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.74"
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "66.227.62"
$ws.Range("E16").Value = "  +5.80%  "
$ws.Range("D17").Value = "3.498.57"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("D19").Value = "3.003.75"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.51"
$ws.Range("E20").Value = "  +5.61%  "
$ws.Range("E21").Value = "  +5.44%  "
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("E23").Value = "  +6.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.48"
$ws.Range("E24").Value = "  +4.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +14.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.29"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.43"
$ws.Range("E27").Value = "  +4.68%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +14.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("E30").Value = "  +18.11%  "
$ws.Range("E31").Value = "  +5.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000104"
$ws.Range("E32").Value = "  -5.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.38"
$ws.Range("E33").Value = "  +5.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +4.04%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.84"
$ws.Range("E37").Value = "  +7.79%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.08"
$ws.Range("E38").Value = "  +9.25%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.82"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("E41").Value = "  +16.36%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.09"
$ws.Range("E42").Value = "  +7.14%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.122"
$ws.Range("E43").Value = "  +6.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.45"
$ws.Range("E44").Value = "  +4.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "404.19"
$ws.Range("E45").Value = "  +13.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0358"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").Value = "2.791.44"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.93"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.93"
$ws.Range("E50").Value = "  +11.14%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.16"
$ws.Range("E51").Value = "  +10.84%  "
